$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 296, shifting existing rows (296-346) down to (297-347)
$ws.Rows.Item(296).Insert()

# Populate the newly inserted row 296 with the new weekly data point
$ws.Cells.Item(296, 1).Value2 = 7
$ws.Cells.Item(296, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(296, 3).Value2 = "Ñuble"
$ws.Cells.Item(296, 4).Value2 = 45131
$ws.Cells.Item(296, 5).Value2 = 16
$ws.Cells.Item(296, 6).Value2 = "Fruta"
$ws.Cells.Item(296, 7).Value2 = 100101
$ws.Cells.Item(296, 8).Value2 = "Berries"
$ws.Cells.Item(296, 9).Value2 = 100101007
$ws.Cells.Item(296, 10).Value2 = "Kiwi"
$ws.Cells.Item(296, 11).Value2 = "Hayward"
$ws.Cells.Item(296, 12).Value2 = "Primera"
$ws.Cells.Item(296, 13).Value2 = 120
$ws.Cells.Item(296, 14).Value2 = 12000
$ws.Cells.Item(296, 15).Value2 = 13000
$ws.Cells.Item(296, 16).Value2 = 12500
$ws.Cells.Item(296, 17).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(296, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(296, 19).Value2 = 694
$ws.Cells.Item(296, 20).Value2 = 18

Write-Host "Row inserted and populated"
